$d = $word.ActiveDocument

# The translated / final replacement text for the "Perseus observation periods" line.
# In the source document this sentence is built out of many separately-formatted
# runs (leftover "2018" observation dates). The commit replaces the whole run
# sequence inside each of these paragraphs with a single, unformatted run
# containing the newly translated dates.
$newText = "havainnointijaksot vuonna Perseus: 16.-25. Tammikuuta, 7.-16."

$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# Walk backwards so replacing a paragraph's contents never invalidates the
# index of a paragraph we still need to visit.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "Perseus*havainnointijaksot vuonna 2018*") {
        $r = $p.Range
        # Exclude the trailing paragraph mark so only the run content
        # (i.e. everything between <w:pPr> and </w:p>) gets replaced;
        # the paragraph's own identity/properties are left untouched.
        $body = $d.Range($r.Start, $r.End - 1)
        $body.InsertXML($newXml)
    }
}
